$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column H, matching style of existing header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

$values = @(0,0,1,0,1,0,0,0,0,0,0,1,0,0,1,0,0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}

Write-Host "Done"
